# Regenerate save_data "K" column (Strike# -> K) values on Sheet1.
# Only column G (header "K") changes for this workbook; rows 3 and 9 keep
# their existing values, all other listed rows get new computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    10 = 2
    11 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
